$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before H, shifting old H:K to I:L
$ws.Columns("H:H").Insert()

# Set new column H header/value
$ws.Range("H1").Value = "cuenta"
$ws.Range("H1").Style = $ws.Range("G1").Style
$ws.Columns("H:H").ColumnWidth = 16.36328125

$ws.Range("H2").Value = "CTS CLIENTES"

# Update shifted columns' row2 values to new data
$ws.Range("I2").Value = "FAILED"
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "30 jun. 2023, 14:13:39"
$ws.Range("L2").Value = ""

# Update conditional formatting range from J6 to K5
$ws.Cells.FormatConditions.Delete()
